# Update the recomputed "Return_with_prediction" (G), "return_pct_change" (H)
# and "mean_return_pct_change" (I, row 2 only) columns with refreshed values
# for rows 2-57 of Sheet1, as produced by the latest recurrence run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1071440481298471
$ws.Range("H2").Value = 13.11554920545
$ws.Range("I2").Value = 10.96991618690094
$ws.Range("G3").Value = 0.09749306358327187
$ws.Range("H3").Value = 47.50688679324799
$ws.Range("G4").Value = 0.01034666405895339
$ws.Range("H4").Value = 10.36126695246682
$ws.Range("G5").Value = 0.04456145850098477
$ws.Range("H5").Value = 432.8012175552805
$ws.Range("G6").Value = -0.2241164520860654
$ws.Range("H6").Value = -1.320443497279872
$ws.Range("G7").Value = -0.2084579080731468
$ws.Range("H7").Value = 16.57535614129695
$ws.Range("G8").Value = -0.3561977129110726
$ws.Range("H8").Value = 3.796876157055699
$ws.Range("G9").Value = -0.3814893162140143
$ws.Range("H9").Value = 4.300947123901501
$ws.Range("G10").Value = -0.04162910228201853
$ws.Range("H10").Value = -356.8922512480014
$ws.Range("G11").Value = 0.102816102469102
$ws.Range("H11").Value = 739.966086841599
$ws.Range("G12").Value = 0.2387273641831133
$ws.Range("H12").Value = 5.078843261599995
$ws.Range("G13").Value = 0.3245399048601543
$ws.Range("H13").Value = 23.23960808688511
$ws.Range("G14").Value = -0.03503354394610106
$ws.Range("H14").Value = -266.181106952137
$ws.Range("G15").Value = 0.0144250802314894
$ws.Range("H15").Value = -28.54172239384309
$ws.Range("G16").Value = 0.1412727757540835
$ws.Range("H16").Value = 19.72359901731749
$ws.Range("G17").Value = 0.1947365578180325
$ws.Range("H17").Value = -11.01130686400912
$ws.Range("G18").Value = 0.04517901313094688
$ws.Range("H18").Value = -25.2827006460474
$ws.Range("G19").Value = 0.05748703044335348
$ws.Range("H19").Value = -36.1888394936304
$ws.Range("G20").Value = -0.1482516631290654
$ws.Range("H20").Value = -1.86949080245939
$ws.Range("G21").Value = -0.140465427323711
$ws.Range("H21").Value = 29.70750038318037
$ws.Range("G22").Value = 0.04791281145783299
$ws.Range("H22").Value = -11.90593925146042
$ws.Range("G23").Value = 0.03705856480641802
$ws.Range("H23").Value = -9.259468514868953
$ws.Range("G24").Value = 0.1387876099029383
$ws.Range("H24").Value = 19.91914433563659
$ws.Range("G25").Value = 0.1495183491146239
$ws.Range("H25").Value = -1.681842617198704
$ws.Range("G26").Value = -0.0167802859385372
$ws.Range("H26").Value = -131.7346929703294
$ws.Range("G27").Value = 0.00713120981305596
$ws.Range("H27").Value = -85.86950479180057
$ws.Range("G28").Value = 0.1546322568177046
$ws.Range("H28").Value = 1.12561133181781
$ws.Range("G29").Value = 0.1829370678587904
$ws.Range("H29").Value = 7.166818170791295
$ws.Range("G30").Value = -0.0000300124910226448
$ws.Range("H30").Value = -100.1533903596223
$ws.Range("G31").Value = 0.0526065813519051
$ws.Range("H31").Value = 442.0581216384637
$ws.Range("G32").Value = 0.0274864649565908
$ws.Range("H32").Value = -26.29645171145957
$ws.Range("G33").Value = 0.01255556853970396
$ws.Range("H33").Value = -51.90107937989426
$ws.Range("G34").Value = 0.1088406164671586
$ws.Range("H34").Value = -14.94701087154319
$ws.Range("G35").Value = 0.1506219526517972
$ws.Range("H35").Value = 17.06921993648769
$ws.Range("G36").Value = -0.02386636681214747
$ws.Range("H36").Value = -258.7715905246091
$ws.Range("G37").Value = -0.06969768590559097
$ws.Range("H37").Value = -555.1094628328675
$ws.Range("G38").Value = -0.02613837875793051
$ws.Range("H38").Value = -1179.778929594417
$ws.Range("G39").Value = -0.009700677248651304
$ws.Range("H39").Value = 70.96406484175579
$ws.Range("G40").Value = 0.1142724963390129
$ws.Range("H40").Value = -22.55284517019041
$ws.Range("G41").Value = 0.1125492283213453
$ws.Range("H41").Value = -30.26655394667407
$ws.Range("G42").Value = 0.03782356034534425
$ws.Range("H42").Value = -41.41777051253108
$ws.Range("G43").Value = 0.06987325344147276
$ws.Range("H43").Value = 101.0132386471804
$ws.Range("G44").Value = 0.02113111310401822
$ws.Range("H44").Value = 49.73030903377693
$ws.Range("G45").Value = 0.01823040605181612
$ws.Range("H45").Value = -55.59844862316174
$ws.Range("G46").Value = -0.03897816130772801
$ws.Range("H46").Value = 40.78055022985077
$ws.Range("G47").Value = -0.02149700464790572
$ws.Range("H47").Value = 47.96136697828252
$ws.Range("G48").Value = -0.1447771925066553
$ws.Range("H48").Value = -14.92496287165941
$ws.Range("G49").Value = -0.1528470071380833
$ws.Range("H49").Value = 22.60182374610947
$ws.Range("G50").Value = 0.1053319682125969
$ws.Range("H50").Value = -3.253018186508006
$ws.Range("G51").Value = 0.1747138661350466
$ws.Range("H51").Value = 74.24196951510076
$ws.Range("G52").Value = 0.05643414876138974
$ws.Range("H52").Value = -5.342091718695986
$ws.Range("G53").Value = 0.075264696479983
$ws.Range("H53").Value = 11.43279476744526
$ws.Range("G54").Value = -0.0448042233425162
$ws.Range("H54").Value = 35.92126505362204
$ws.Range("G55").Value = -0.03429888082898155
$ws.Range("H55").Value = 55.59222270591358
$ws.Range("G56").Value = 0.0822884497442096
$ws.Range("H56").Value = 79.56687817256505
$ws.Range("G57").Value = 0.08370933566470429
$ws.Range("H57").Value = 1519.049086189271
